# Fruta / hortaliza, semanal
# Insert a new weekly record for "Membrillo" (Vega Modelo de Temuco) at row 87,
# shifting all following records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 87; rows 87-153 shift down to 88-154
$ws.Rows.Item(87).Insert()

# Fill the newly inserted row 87 with the new record's data
$ws.Range("A87").Value = 10
$ws.Range("B87").Value = "Vega Modelo de Temuco"
$ws.Range("C87").Value = "La Araucanía"
$ws.Range("D87").Value = 44669
$ws.Range("E87").Value = 9
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100104
$ws.Range("H87").Value = "Frutos de pepita"
$ws.Range("I87").Value = 100104003
$ws.Range("J87").Value = "Membrillo"
$ws.Range("K87").Value = "Champion"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 70
$ws.Range("N87").Value = 14000
$ws.Range("O87").Value = 15000
$ws.Range("P87").Value = 14571
$ws.Range("Q87").Value = "$/bandeja 18 kilos granel"
$ws.Range("R87").Value = "Región de O'Higgins"
$ws.Range("S87").Value = 810
$ws.Range("T87").Value = 18
